# Wario 100% - Level 10 begin.
# Insert a new "Checkpoint" row before the existing row 100 ("Use Key"),
# shifting every row from 100 onward down by one, then populate the new
# row's data, add the two newly-recorded times on the rows that used to
# be 100/101/103 (now 101/102/104), and finally restore the view state
# (frozen-pane top-left cell + active selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a blank row at row 100 (pushes old rows 100:163 down to 101:164,
# updates the shared formulas' ranges and the merged-cell addresses).
$ws.Rows.Item(100).Insert()

# Copy the formatting (style id, borders, fill, number format) from the row
# directly below (the shifted former row 100) onto the brand new row so the
# blank row matches the sheet's normal data-row styling.
$ws.Range("A101:D101").Copy()
$ws.Range("A100:D100").PasteSpecial(-4122)

# Fill in the new checkpoint row's data.
$ws.Range("A100").Value = "Checkpoint"
$ws.Range("B100").Value = 40410
$ws.Range("C100").Value = 32946
$ws.Range("D100").Formula = "=IF(B100>0,C100-B100,0)"

# The next two rows (former rows 100 and 101, "Use Key" / "Exit Level") and
# the former row 103 ("Enter Level 10") now record their Pre-Level time too.
$ws.Range("B101").Value = 42018
$ws.Range("B102").Value = 42221
$ws.Range("B104").Value = 42470

# Restore the sheet view: frozen pane scrolled to A94 and the active
# selection on B105 (same relative spot the author was working at).
$ws.Application.ActiveWindow.ScrollRow = 94
$ws.Range("B105").Select()
